$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1) and "全部类型" sheet (sheet4)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 614
    $ws.Range("F3").Value = 3761
}
